# Weekly update: insert two new rows of fresh "Cebolla" price data right
# above the previous newest "1a (guarda)" / "2a (guarda)" O'Higgins rows,
# shifting all the subsequent rows down by two (dimension grows from
# A1:R413 to A1:R415).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows at 361:362, pushing the old rows 361-413 down to
# become rows 363-415.
$ws.Rows("361:362").Insert()

# New row 361: "1a (guarda)" price point for the latest week.
$ws.Range("A361").Value = 8
$ws.Range("B361").Value = "Terminal La Palmera de La Serena"
$ws.Range("C361").Value = "Coquimbo"
$ws.Range("D361").Value = 44474
$ws.Range("E361").Value = 4
$ws.Range("F361").Value = 100112004
$ws.Range("G361").Value = "Cebolla"
$ws.Range("H361").Value = "Sin especificar"
$ws.Range("I361").Value = "1a (guarda)"
$ws.Range("J361").Value = 3000
$ws.Range("K361").Value = 4800
$ws.Range("L361").Value = 5000
$ws.Range("M361").Value = 4900
$ws.Range("N361").Value = "`$/malla 16 kilos"
$ws.Range("O361").Value = "Región de O'Higgins"
$ws.Range("P361").Value = 306
$ws.Range("Q361").Value = 16
$ws.Range("R361").Value = "Hortaliza"

# New row 362: "2a (guarda)" price point for the latest week.
$ws.Range("A362").Value = 8
$ws.Range("B362").Value = "Terminal La Palmera de La Serena"
$ws.Range("C362").Value = "Coquimbo"
$ws.Range("D362").Value = 44474
$ws.Range("E362").Value = 4
$ws.Range("F362").Value = 100112004
$ws.Range("G362").Value = "Cebolla"
$ws.Range("H362").Value = "Sin especificar"
$ws.Range("I362").Value = "2a (guarda)"
$ws.Range("J362").Value = 1600
$ws.Range("K362").Value = 4500
$ws.Range("L362").Value = 4600
$ws.Range("M362").Value = 4550
$ws.Range("N362").Value = "`$/malla 16 kilos"
$ws.Range("O362").Value = "Región de O'Higgins"
$ws.Range("P362").Value = 284
$ws.Range("Q362").Value = 16
$ws.Range("R362").Value = "Hortaliza"
